$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 3542.88
$ws.Range("I33").Value = 4710.1177
$ws.Range("K33").Value = 4710.1177
$ws.Range("M33").Value = -4481.1177
$ws.Range("H51").Value = 52031.188
$ws.Range("J51").Value = 87666.78
$ws.Range("L51").Value = 87666.78
$ws.Range("N51").Value = -88634.78
$ws.Range("H69").Value = 6625
$ws.Range("J69").Value = 6625
$ws.Range("L69").Value = 19875
$ws.Range("N69").Value = -21623
$ws.Range("H72").Value = 6625
$ws.Range("J72").Value = 6625
$ws.Range("L72").Value = 59625
$ws.Range("N72").Value = -68361
$ws.Range("H86").Value = 2450.2942
$ws.Range("I86").Value = 1928
$ws.Range("J86").Value = 2815.9
$ws.Range("K86").Value = 1928
$ws.Range("L86").Value = 2815.9
$ws.Range("M86").Value = -805
$ws.Range("N86").Value = -5061.9
$ws.Range("H89").Value = 2450.2942
$ws.Range("I89").Value = 1928
$ws.Range("J89").Value = 2815.9
$ws.Range("K89").Value = 9640
$ws.Range("L89").Value = 14079.5
$ws.Range("M89").Value = -4024
$ws.Range("N89").Value = -25311.5
$ws.Range("H129").Value = 1487.4
$ws.Range("I129").Value = 980.7692
$ws.Range("J129").Value = 2428.2856
$ws.Range("K129").Value = 2942.3076
$ws.Range("L129").Value = 7284.8568
$ws.Range("M129").Value = 2057.6924
$ws.Range("N129").Value = -17284.8568
$ws.Range("H138").Value = 2581.6667
$ws.Range("I138").Value = 1863.0869
$ws.Range("J138").Value = 3332.9092
$ws.Range("K138").Value = 5589.2607
$ws.Range("L138").Value = 9998.7276
$ws.Range("M138").Value = -449.2606999999998
$ws.Range("N138").Value = -20278.7276

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 50000
$ws.Range("I9").Value = 50000
$ws.Range("K9").Value = 50000
$ws.Range("M9").Value = -49830
$ws.Range("H20").Value = 50000
$ws.Range("I20").Value = 50000
$ws.Range("K20").Value = 50000
$ws.Range("M20").Value = -49730
$ws.Range("H32").Value = 32281.834
$ws.Range("I32").Value = 32281.834
$ws.Range("K32").Value = 32281.834
$ws.Range("M32").Value = -31994.834
$ws.Range("H61").Value = 3284
$ws.Range("I61").Value = 3038.5806
$ws.Range("K61").Value = 3038.5806
$ws.Range("M61").Value = -2826.5806
$ws.Range("H136").Value = 3284
$ws.Range("I136").Value = 3038.5806
$ws.Range("K136").Value = 9115.7418
$ws.Range("M136").Value = -6565.7418

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 43491588
$ws.Range("I105").Value = 55571724
$ws.Range("K105").Value = 55571724
$ws.Range("M105").Value = -55569977
$ws.Range("H107").Value = 1356.4166
$ws.Range("I107").Value = 1212.6666
$ws.Range("K107").Value = 1212.6666
$ws.Range("M107").Value = 707.3334
$ws.Range("H134").Value = 3241.875
$ws.Range("I134").Value = 1842.4073
$ws.Range("K134").Value = 5527.2219
$ws.Range("M134").Value = -2992.2219

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2928.8965
$ws.Range("I58").Value = 1918.25
$ws.Range("J58").Value = 7780
$ws.Range("K58").Value = 1918.25
$ws.Range("L58").Value = 7780
$ws.Range("M58").Value = -1715.25
$ws.Range("N58").Value = -8186
$ws.Range("H132").Value = 2510.5938
$ws.Range("I132").Value = 2096.0356
$ws.Range("J132").Value = 5412.5
$ws.Range("K132").Value = 6288.1068
$ws.Range("L132").Value = 16237.5
$ws.Range("M132").Value = -3758.1068
$ws.Range("N132").Value = -21297.5
$ws.Range("H136").Value = 2928.8965
$ws.Range("I136").Value = 1918.25
$ws.Range("J136").Value = 7780
$ws.Range("K136").Value = 5754.75
$ws.Range("L136").Value = 23340
$ws.Range("M136").Value = -3204.75
$ws.Range("N136").Value = -28440

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 2020909.1
$ws.Range("I2").Value = 3225877.5
$ws.Range("K2").Value = 19355265
$ws.Range("M2").Value = -19355152
$ws.Range("H5").Value = 1337.5
$ws.Range("J5").Value = 1483.3334
$ws.Range("L5").Value = 4450.0002
$ws.Range("N5").Value = -4674.0002
$ws.Range("H21").Value = 537.25
$ws.Range("I21").Value = 483
$ws.Range("K21").Value = 1449
$ws.Range("M21").Value = -1276
$ws.Range("H38").Value = 57.72
$ws.Range("I38").Value = 37.285713
$ws.Range("J38").Value = 83.72727
$ws.Range("K38").Value = 111.857139
$ws.Range("L38").Value = 251.18181
$ws.Range("M38").Value = 235.142861
$ws.Range("N38").Value = -945.18181
$ws.Range("H39").Value = 3992.7407
$ws.Range("J39").Value = 4800
$ws.Range("L39").Value = 14400
$ws.Range("N39").Value = -14988
$ws.Range("H55").Value = 1205.9333
$ws.Range("I55").Value = 230.83333
$ws.Range("J55").Value = 1856
$ws.Range("K55").Value = 692.49999
$ws.Range("L55").Value = 5568
$ws.Range("M55").Value = -515.49999
$ws.Range("N55").Value = -5922
$ws.Range("H117").Value = 238274
$ws.Range("J117").Value = 416812.25
$ws.Range("L117").Value = 1250436.75
$ws.Range("N117").Value = -1257320.75
$ws.Range("H122").Value = 23473.334
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 3469.4
$ws.Range("J132").Value = 2415.3333
$ws.Range("L132").Value = 21737.9997
$ws.Range("N132").Value = -26797.9997
$ws.Range("H135").Value = 1337.5
$ws.Range("J135").Value = 1483.3334
$ws.Range("L135").Value = 13350.0006
$ws.Range("N135").Value = -18420.0006

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 1121.2
$ws.Range("I29").Value = 1121.2
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 1121.2
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -831.2
$ws.Range("N29").ClearContents()
$ws.Range("H41").Value = 2704.4
$ws.Range("I41").Value = 1130.5
$ws.Range("K41").Value = 1130.5
$ws.Range("M41").Value = -775.5
$ws.Range("H97").Value = 922.3333
$ws.Range("I97").Value = 1183.25
$ws.Range("J97").Value = 520.9231
$ws.Range("K97").Value = 1183.25
$ws.Range("L97").Value = 520.9231
$ws.Range("M97").Value = -687.25
$ws.Range("N97").Value = -1512.9231
$ws.Range("H132").Value = 7161.5386
$ws.Range("I132").Value = 5961.905
$ws.Range("J132").Value = 12200
$ws.Range("K132").Value = 17885.715
$ws.Range("L132").Value = 36600
$ws.Range("M132").Value = -15355.715
$ws.Range("N132").Value = -41660
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3208.5386
$ws.Range("I22").Value = 1940.75
$ws.Range("K22").Value = 1940.75
$ws.Range("M22").Value = -1645.75
$ws.Range("H27").Value = 3208.5386
$ws.Range("I27").Value = 1940.75
$ws.Range("K27").Value = 1940.75
$ws.Range("M27").Value = -1833.75
$ws.Range("H122").Value = 8148.25
$ws.Range("I122").Value = 7109
$ws.Range("K122").Value = 21327
$ws.Range("M122").Value = -18877
$ws.Range("H132").Value = 5399.488
$ws.Range("I132").Value = 3200.1333
$ws.Range("J132").Value = 11397.728
$ws.Range("K132").Value = 9600.3999
$ws.Range("L132").Value = 34193.18399999999
$ws.Range("M132").Value = -7070.3999
$ws.Range("N132").Value = -39253.18399999999
$ws.Range("H136").Value = 3348.2
$ws.Range("I136").Value = 2745.3704
$ws.Range("J136").Value = 4600.231
$ws.Range("K136").Value = 8236.111199999999
$ws.Range("L136").Value = 13800.693
$ws.Range("M136").Value = -5686.111199999999
$ws.Range("N136").Value = -18900.693

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 75000
$ws.Range("J64").Value = 75000
$ws.Range("L64").Value = 75000
$ws.Range("N64").Value = -75496
$ws.Range("H67").Value = 75000
$ws.Range("J67").Value = 75000
$ws.Range("L67").Value = 75000
$ws.Range("N67").Value = -76716
$ws.Range("H95").Value = 67772
$ws.Range("J95").Value = 67772
$ws.Range("L95").Value = 67772
$ws.Range("N95").Value = -73264
$ws.Range("H113").Value = 419.63635
$ws.Range("J113").Value = 627.3333
$ws.Range("L113").Value = 1881.9999
$ws.Range("N113").Value = -6221.9999
$ws.Range("H124").Value = 106993.86
$ws.Range("J124").Value = 106993.86
$ws.Range("L124").Value = 106993.86
$ws.Range("N124").Value = -116813.86
$ws.Range("H132").Value = 4496.3022
$ws.Range("I132").Value = 2362.926
$ws.Range("K132").Value = 7088.778
$ws.Range("M132").Value = -4558.778
$ws.Range("H136").Value = 2190.7715
$ws.Range("I136").Value = 1450.3043
$ws.Range("K136").Value = 4350.9129
$ws.Range("M136").Value = -1800.9129
